# Applies the cryptocurrency price/volume refresh described in the commit.
# All Price (D) and Volume(1h) (E) values are stored as text in this sheet,
# so numeric-looking Price values are written with a leading apostrophe and
# the cell style is reset to "Normal" afterwards to avoid Excel silently
# converting them to numbers / changing their number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.328.76"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.843.50"
$ws.Range("D4").Value = "'0.9964"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'239.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'0.6262"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.9985"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "'0.2897"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'0.07725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.843.98"
$ws.Range("D13").Value = "'4.986"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'0.6803"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "'0.00001049"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "'6.181"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "29.385.85"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'228.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'0.9982"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'7.478"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'0.9978"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'158.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'8.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'0.1369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "'17.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "'0.06411"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.55%  "
$ws.Range("D29").Value = "'1.400"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "'1.480"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "'4.087"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "'0.6971"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "'2.579"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "1.268.84"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").Value = "'2.839"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.44%  "
$ws.Range("D39").Value = "'0.01834"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").Value = "'6.733"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("D41").Value = "'0.9146"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'0.9978"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "2.006.59"
$ws.Range("E43").Value = "  -18.46%  "
$ws.Range("D44").Value = "'101.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'66.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'7.077"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").Value = "'1.724"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "'0.1163"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.3962"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05694"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
